# Updated cryptos list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns for each coin row. Values that look like plain
# numbers are prefixed with a leading apostrophe so Excel stores them as
# text (matching the original inline-string "text" cells) instead of
# silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.406.73"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.068.63"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'235.34"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'57.52"
$ws.Range("E8").Value = "  -0.60%  "
$ws.Range("D9").Value = "'0.396"
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("D10").Value = "'0.0772"
$ws.Range("D11").Value = "'0.102"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "2.373.07"
$ws.Range("D13").Value = "'14.46"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "'20.78"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("D15").Value = "'0.782"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "'5.20"
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "2.069.87"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").Value = "37.359.44"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'6.35"
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "'69.74"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'226.20"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'167.22"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("D27").Value = "'8.87"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -4.15%  "
$ws.Range("D29").Value = "'19.10"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'4.53"
$ws.Range("E32").Value = "  +0.65%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "'4.57"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").Value = "'2.50"
$ws.Range("E35").Value = "  -2.87%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").Value = "'5.68"
$ws.Range("D40").Value = "'0.0968"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'98.31"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("D43").Value = "1.481.65"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("E45").Value = "  +0.65%  "
$ws.Range("D46").Value = "'4.09"
$ws.Range("E46").Value = "  -10.09%  "
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("D48").Value = "'15.36"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("D51").Value = "2.259.88"
$ws.Range("E51").Value = "  +0.31%  "
